$wb = $excel.ActiveWorkbook

# "DATA" sheet holds the browser/version test matrix. Swap the old,
# now-unused Chrome/Edge browser versions for the Selenium-grid-via-docker
# setup: row 4 (chrome) gets version "98.0", row 3's browser becomes
# "firefox" with version "97.0".
$ws = $wb.Worksheets.Item("DATA")

# Leading apostrophe forces these to stay text cells (matching the sheet's
# original text-typed version column) instead of being auto-coerced to
# numbers by Excel's smart input parsing.
$ws.Range("D4").Value = "'98.0"
$ws.Range("C3").Value = "firefox"
$ws.Range("D3").Value = "'97.0"
